$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 244, pushing existing rows 244:268 down to 245:269
$ws.Rows.Item(244).Insert()

# Populate the newly inserted row 244 with the new price record.
# (Columns A,B,C,E,F,G,H,I,J,K,R mirror the constant values used throughout
# this block of "Macroferia Regional de Talca - Piña / Caramelo" rows.)
$ws.Cells.Item(244, 1).Value = 5
$ws.Cells.Item(244, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(244, 3).Value = "Maule"
$ws.Cells.Item(244, 4).Value = 44769
$ws.Cells.Item(244, 5).Value = 7
$ws.Cells.Item(244, 6).Value = "Fruta"
$ws.Cells.Item(244, 7).Value = 100108
$ws.Cells.Item(244, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(244, 9).Value = 100108005
$ws.Cells.Item(244, 10).Value = "Piña"
$ws.Cells.Item(244, 11).Value = "Caramelo"
$ws.Cells.Item(244, 12).Value = "Segunda"
$ws.Cells.Item(244, 13).Value = 180
$ws.Cells.Item(244, 14).Value = 19000
$ws.Cells.Item(244, 15).Value = 19000
$ws.Cells.Item(244, 16).Value = 19000
$ws.Cells.Item(244, 17).Value = "$/caja 14 unidades"
$ws.Cells.Item(244, 18).Value = "Ecuador"
$ws.Cells.Item(244, 19).Value = 1357
$ws.Cells.Item(244, 20).Value = 14
